# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.368.37"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "3.520.76"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "612.45"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "151.28"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").Value = "3.520.54"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").Value = "7.08"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "0.0000220"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.117.92"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "31.99"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "3.516.17"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "67.379.42"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "6.38"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "15.21"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "444.59"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "9.40"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").Value = "77.36"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "0.0000129"
$ws.Range("E25").Value = "  +8.29%  "
$ws.Range("D26").Value = "3.662.66"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -7.59%  "
$ws.Range("E33").Value = "  +4.19%  "
$ws.Range("D34").Value = "25.86"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "3.513.63"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("D38").Value = "7.99"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "177.14"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("D43").Value = "0.0879"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("D45").Value = "0.881"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "28.43"
$ws.Range("E46").Value = "  -3.62%  "
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("D50").Value = "7.58"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "0.993"
$ws.Range("E51").Value = "  -2.00%  "
